# Update automàtic: dades i banners [2026-02-15 23:50]
# Applies the scraped MeteoCat refresh: new extraction timestamps plus
# updated observation values for the affected stations/columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/number-like values (dates, "X.X mm", "X.X °C", "X hPa", ...) ---
# These do not get reinterpreted by Excel, so a direct .Value assignment is safe.
$ws.Range("E2").Value = "2026-02-15 23:48:26"
$ws.Range("E3").Value = "2026-02-15 23:48:28"
$ws.Range("I3").Value = "3.3 mm"
$ws.Range("O3").Value = "-4.6 °C"
$ws.Range("E4").Value = "2026-02-15 23:48:31"
$ws.Range("O4").Value = "7.2 °C"
$ws.Range("E5").Value = "2026-02-15 23:48:33"
$ws.Range("I5").Value = "10.0 mm"
$ws.Range("E6").Value = "2026-02-15 23:48:36"
$ws.Range("E7").Value = "2026-02-15 23:48:38"
$ws.Range("E8").Value = "2026-02-15 23:48:41"
$ws.Range("E9").Value = "2026-02-15 23:48:43"
$ws.Range("N9").Value = "5.6 °C 23:28 TU"
$ws.Range("O9").Value = "10.4 °C"
$ws.Range("E10").Value = "2026-02-15 23:48:45"
$ws.Range("E11").Value = "2026-02-15 23:48:46"
$ws.Range("O11").Value = "6.6 °C"
$ws.Range("E12").Value = "2026-02-15 23:48:47"
$ws.Range("N12").Value = "6.1 °C 23:29 TU"
$ws.Range("O12").Value = "10.5 °C"
$ws.Range("E13").Value = "2026-02-15 23:48:48"
$ws.Range("K13").Value = "6.7 MJ/m2"
$ws.Range("E14").Value = "2026-02-15 23:48:49"
$ws.Range("E15").Value = "2026-02-15 23:48:51"
$ws.Range("O15").Value = "10.2 °C"
$ws.Range("E16").Value = "2026-02-15 23:48:52"
$ws.Range("E17").Value = "2026-02-15 23:48:53"
$ws.Range("O17").Value = "3.3 °C"
$ws.Range("E18").Value = "2026-02-15 23:48:54"
$ws.Range("O18").Value = "7.4 °C"
$ws.Range("E19").Value = "2026-02-15 23:48:55"
$ws.Range("O19").Value = "3.7 °C"
$ws.Range("E20").Value = "2026-02-15 23:48:56"
$ws.Range("O20").Value = "-2.4 °C"
$ws.Range("E21").Value = "2026-02-15 23:48:57"
$ws.Range("J21").Value = "1015.2 hPa"
$ws.Range("O21").Value = "7.6 °C"
$ws.Range("E22").Value = "2026-02-15 23:49:00"
$ws.Range("N22").Value = "-6.5 °C 23:10 TU"
$ws.Range("E23").Value = "2026-02-15 23:49:02"
$ws.Range("I23").Value = "6.8 mm"
$ws.Range("E24").Value = "2026-02-15 23:49:04"
$ws.Range("E25").Value = "2026-02-15 23:49:07"
$ws.Range("O25").Value = "-1.1 °C"
$ws.Range("E26").Value = "2026-02-15 23:49:09"
$ws.Range("E27").Value = "2026-02-15 23:49:11"
$ws.Range("O27").Value = "0.2 °C"
$ws.Range("E28").Value = "2026-02-15 23:49:14"
$ws.Range("E29").Value = "2026-02-15 23:49:16"
$ws.Range("O29").Value = "9.8 °C"
$ws.Range("E30").Value = "2026-02-15 23:49:18"
$ws.Range("O30").Value = "9.6 °C"
$ws.Range("E31").Value = "2026-02-15 23:49:21"
$ws.Range("O31").Value = "10.4 °C"
$ws.Range("E32").Value = "2026-02-15 23:49:23"
$ws.Range("E33").Value = "2026-02-15 23:49:25"
$ws.Range("O33").Value = "6.2 °C"
$ws.Range("E34").Value = "2026-02-15 23:49:27"
$ws.Range("E35").Value = "2026-02-15 23:49:30"
$ws.Range("J35").Value = "1019.6 hPa"
$ws.Range("E36").Value = "2026-02-15 23:49:33"
$ws.Range("O36").Value = "11.1 °C"
$ws.Range("E37").Value = "2026-02-15 23:49:35"
$ws.Range("J37").Value = "1016.7 hPa"
$ws.Range("O37").Value = "5.6 °C"
$ws.Range("E38").Value = "2026-02-15 23:49:37"
$ws.Range("E39").Value = "2026-02-15 23:49:40"
$ws.Range("O39").Value = "-2.4 °C"
$ws.Range("E40").Value = "2026-02-15 23:49:42"
$ws.Range("J40").Value = "1016.5 hPa"
$ws.Range("O40").Value = "8.2 °C"
$ws.Range("E41").Value = "2026-02-15 23:49:45"
$ws.Range("E42").Value = "2026-02-15 23:49:47"
$ws.Range("O42").Value = "10.3 °C"
$ws.Range("E43").Value = "2026-02-15 23:49:50"
$ws.Range("O43").Value = "6.4 °C"
$ws.Range("E44").Value = "2026-02-15 23:49:52"
$ws.Range("I44").Value = "7.3 mm"
$ws.Range("E45").Value = "2026-02-15 23:49:55"
$ws.Range("I45").Value = "5.9 mm"
$ws.Range("J45").Value = "1023.1 hPa"
$ws.Range("E46").Value = "2026-02-15 23:49:57"
$ws.Range("J46").Value = "1019.5 hPa"

# --- Humidity percentages ---
# Excel auto-parses bare "NN%" text into a numeric percentage, which would
# silently change both the stored value and the cell style. Force the cells
# to Text format first, assign, then restore the original (General) number
# format/style by pasting formats from an untouched same-column cell (H2),
# so the resulting cell keeps its original style with a plain text value.
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "74%"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "63%"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "56%"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "50%"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "62%"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "42%"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "75%"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "65%"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "43%"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "62%"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "81%"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "53%"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "45%"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "80%"

$ws.Range("H2").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("H36").PasteSpecial(-4122)
$ws.Range("H40").PasteSpecial(-4122)
$ws.Range("H44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

